# Updates the cryptos price/volume table (and a handful of re-ranked
# coin name/link/price rows) to the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value.
$updates = @{
    "D2" = "27.753.28"
    "E2" = "  -1.95%  "
    "D3" = "1.896.75"
    "E3" = "  -1.57%  "
    "D4" = "0.9995"
    "E4" = "  -0.62%  "
    "D5" = "312.37"
    "E5" = "  -1.42%  "
    "D6" = "0.9991"
    "E6" = "  -0.57%  "
    "D7" = "0.4905"
    "E7" = "  +1.25%  "
    "D8" = "0.3804"
    "E8" = "  -1.50%  "
    "E9" = "  -1.18%  "
    "D10" = "0.9076"
    "E10" = "  -4.07%  "
    "D11" = "20.86"
    "E11" = "  -0.63%  "
    "D12" = "0.07614"
    "E12" = "  -2.68%  "
    "D13" = "1.871.43"
    "E13" = "  -2.76%  "
    "D14" = "5.478"
    "E14" = "  -1.12%  "
    "D15" = "6.637"
    "E15" = "  -0.63%  "
    "D16" = "90.92"
    "E16" = "  -0.97%  "
    "D17" = "1.000"
    "E17" = "  -0.62%  "
    "D18" = "0.000008720"
    "E18" = "  -1.93%  "
    "D19" = "0.9993"
    "E19" = "  -0.57%  "
    "D20" = "27.747.00"
    "E20" = "  -2.01%  "
    "D21" = "14.53"
    "E21" = "  -2.95%  "
    "D22" = "5.142"
    "E22" = "  -0.89%  "
    "B23" = "Cosmos"
    "C23" = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
    "D23" = "10.75"
    "E23" = "  -2.22%  "
    "B24" = "Monero"
    "C24" = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
    "D24" = "154.06"
    "E24" = "  -1.46%  "
    "B25" = "Toncoin"
    "C25" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "D25" = "1.858"
    "E25" = "  -3.92%  "
    "B26" = "LidoDAOToken"
    "C26" = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
    "D26" = "2.199"
    "E26" = "  +4.01%  "
    "B27" = "EthereumClassic"
    "C27" = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
    "D27" = "18.37"
    "E27" = "  -1.51%  "
    "B28" = "BitcoinCash"
    "C28" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "D28" = "114.96"
    "E28" = "  -1.72%  "
    "B29" = "InternetComputer(DFINITY)"
    "C29" = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
    "D29" = "4.899"
    "E29" = "  -2.56%  "
    "B30" = "Stellar"
    "C30" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D30" = "0.08920"
    "E30" = "  -0.04%  "
    "B31" = "HuobiToken"
    "C31" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D31" = "3.226"
    "E31" = "  -3.96%  "
    "B32" = "ARBITRUM"
    "C32" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D32" = "1.237"
    "E32" = "  -2.75%  "
    "B33" = "ImmutableX"
    "C33" = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
    "D33" = "0.7712"
    "E33" = "  -0.94%  "
    "B34" = "Filecoin"
    "C34" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D34" = "4.643"
    "E34" = "  -1.81%  "
    "B35" = "VeChain"
    "C35" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "D35" = "0.02057"
    "E35" = "  -0.41%  "
    "B36" = "RenderToken"
    "C36" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D36" = "2.545"
    "E36" = "  -6.80%  "
    "B37" = "TrustWalletToken"
    "C37" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D37" = "1.096"
    "E37" = "  -1.51%  "
    "B38" = "TheSandbox"
    "C38" = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
    "D38" = "0.5507"
    "E38" = "  -1.78%  "
    "B39" = "MXToken"
    "C39" = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D39" = "3.000"
    "E39" = "  -0.60%  "
    "D40" = "0.05263"
    "E40" = "  -2.07%  "
    "B41" = "FraxShare"
    "C41" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D41" = "6.937"
    "E41" = "  -2.17%  "
    "B42" = "Algorand"
    "C42" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D42" = "0.1516"
    "E42" = "  -1.80%  "
    "B43" = "Aptos"
    "C43" = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
    "D43" = "8.435"
    "E43" = "  -2.33%  "
    "B44" = "Quant"
    "C44" = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
    "D44" = "112.07"
    "E44" = "  +4.42%  "
    "B45" = "EnergySwap"
    "C45" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D45" = "10.62"
    "E45" = "  -1.60%  "
    "B46" = "Decentraland"
    "C46" = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
    "D46" = "0.4787"
    "E46" = "  -2.75%  "
    "B47" = "PaxDollar"
    "C47" = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
    "D47" = "0.9986"
    "E47" = "  -0.67%  "
    "B48" = "NEARProtocol"
    "C48" = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
    "D48" = "1.636"
    "E48" = "  -2.33%  "
    "B49" = "Aave"
    "C49" = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
    "D49" = "67.36"
    "E49" = "  -3.11%  "
    "B50" = "Cronos"
    "C50" = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
    "D50" = "0.06052"
    "E50" = "  -1.58%  "
    "B51" = "EOS"
    "C51" = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
    "D51" = "0.8970"
    "E51" = "  -1.89%  "
}

foreach ($cellRef in $updates.Keys) {
    $rng = $ws.Range($cellRef)
    # Force text storage so numeric-looking strings (e.g. "20.86",
    # "0.9995") are not auto-converted to Excel numbers, matching the
    # original inline-string cell type. ClearFormats() afterwards drops
    # the temporary "@" number format so untouched cell styling (s attr)
    # is not left behind.
    $rng.NumberFormat = "@"
    $rng.Value = $updates[$cellRef]
    $rng.ClearFormats()
}
